$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "programa especifico" specify-fields block: rows 137-139 ---
# Duplicate the last existing row (136) three times via a real row
# insert+shift (xlShiftDown). Unlike PasteSpecial/Font-object edits, this
# preserves font-only cell styles (e.g. B136's s="7") verbatim instead of
# collapsing them back to the default style.
$ws.Rows.Item(136).Copy()
$ws.Rows.Item(137).Insert(-4121)
$ws.Rows.Item(136).Copy()
$ws.Rows.Item(138).Insert(-4121)
$ws.Rows.Item(136).Copy()
$ws.Rows.Item(139).Insert(-4121)
$excel.CutCopyMode = 0

# Column A of these rows should look like the other orange "specify"
# headers (A134/A135: Times New Roman, bordered, filled) but with a white
# fill instead of orange. Copy that format, then just swap the fill color.
$srcA = $ws.Range("A134")
$dstA = $ws.Range("A137:A139")
$srcA.Copy()
$dstA.PasteSpecial(-4122)
$dstA.Interior.Color = 16777215
$excel.CutCopyMode = 0

$ws.Range("A137").Value = "QEPE_DGE_SQE_B2_P15_1_participa_programa_Especifique1"
$ws.Range("A138").Value = "QEPE_DGE_SQE_B2_P15_2_participa_programa_Especifique2"
$ws.Range("A139").Value = "QEPE_DGE_SQE_B2_P15_3_participa_programa_Especifique3"

# Column B keeps the row-136 style (orange) carried over by the row
# insert; just swap in the real relevance condition text.
$cond = '${(prefixo)_DGE_SQE_B2_P15_0_escola_participa_programa_especifico}=1'
$ws.Range("B137").Value = $cond
$ws.Range("B138").Value = $cond
$ws.Range("B139").Value = $cond

# Match the sheet's standard row height (15pt, explicit custom-height flag).
$ws.Rows.Item(137).RowHeight = 15
$ws.Rows.Item(138).RowHeight = 15
$ws.Rows.Item(139).RowHeight = 15

# Leave the selection where Excel would after typing the last entry and
# moving one row further down.
$ws.Range("B141").Select()
